# Combine the Defense header row with a new Offense header row.
# 1) Insert a new column at the far left (shifts existing Def headers from A:AE to B:AF)
# 2) Label A1 = "DEFENSE", A2 = "OFFENSE"
# 3) Fill row 2 (B2:AF2) with the matching "_Off" column headers

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing header row one column to the right.
$ws.Columns("A:A").Insert()

# Top-left corner labels.
$ws.Range("A1").Value = "DEFENSE"
$ws.Range("A2").Value = "OFFENSE"

# New Offense header labels, aligned under the existing Defense headers
# (column B corresponds to column A pre-insert, i.e. Year_Def -> Year_Off, etc.)
$offHeaders = @(
    "Year_Off",
    "Rk_Off",
    "Tm_Off",
    "Win_Off",
    "Loss_Off",
    "G_Off",
    "PF_Off",
    "Total Yds_Off",
    "Plays_Off",
    "Y/P_Off",
    "TO_Off",
    "FL_Off",
    "1stD_Off",
    "Cmp_Off",
    "Pass Att_Off",
    "Pass Yds_Off",
    "Pass TD_Off",
    "Int_Off",
    "Pass NY/A_Off",
    "Pass 1stD_Off",
    "Rush Att_Off",
    "Rush Yds_Off",
    "Rush TD_Off",
    "Rush Y/A_Off",
    "Rush 1stD_Off",
    "Pen_Off",
    "Pen Yds_Off",
    "1stPy_Off",
    "Sc%_Off",
    "TO%_Off",
    "EXP_Off"
)

for ($i = 0; $i -lt $offHeaders.Length; $i++) {
    $col = $i + 2   # column B = 2
    $ws.Cells.Item(2, $col).Value = $offHeaders[$i]
}
